# Updates the "cryptos" sheet with refreshed price/volume figures.
# Numeric-looking Price (column D) values are forced to Text before
# assignment (and the style reset to Normal afterwards) so Excel doesn't
# silently coerce strings like "18.09" or "0.789" into floating point
# numbers and lose the original formatting/precision.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.980.10'
$ws.Range('E2').Value = '  +2.79%  '
$ws.Range('D3').Value = '1.600.12'
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('E5').Value = '  +2.44%  '
$ws.Range('E6').Value = '  -0.13%  '
$ws.Range('E7').Value = '  +1.64%  '
$ws.Range('E8').Value = '  +1.67%  '
$ws.Range('E9').Value = '  -0.06%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '18.09'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.97%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0811'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.76%  '
$ws.Range('D12').Value = '1.823.24'
$ws.Range('E12').Value = '  +2.58%  '
$ws.Range('D13').Value = '1.604.67'
$ws.Range('E13').Value = '  +2.72%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.00'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.48%  '
$ws.Range('E15').Value = '  +1.47%  '
$ws.Range('D16').Value = '25.996.83'
$ws.Range('E16').Value = '  +2.81%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '60.20'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.66%  '
$ws.Range('E18').Value = '  +1.50%  '
$ws.Range('E19').Value = '  -0.15%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '201.59'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +8.04%  '
$ws.Range('E21').Value = '  +2.84%  '
$ws.Range('E22').Value = '  +0.04%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.01'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.55%  '
$ws.Range('E24').Value = '  +7.18%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '141.51'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.35%  '
$ws.Range('E26').Value = '  -0.13%  '
$ws.Range('E27').Value = '  -7.16%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.11'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.40%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.45'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.07%  '
$ws.Range('E30').Value = '  +1.73%  '
$ws.Range('E31').Value = '  +2.46%  '
$ws.Range('E32').Value = '  +1.27%  '
$ws.Range('E33').Value = '  -1.09%  '
$ws.Range('E34').Value = '  +0.06%  '
$ws.Range('E35').Value = '  +2.44%  '
$ws.Range('D36').Value = '1.123.34'
$ws.Range('E36').Value = '  +3.21%  '
$ws.Range('E37').Value = '  +10.43%  '
$ws.Range('E38').Value = '  +0.24%  '
$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.32'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.40%  '
$ws.Range('B40').Value = 'ARBITRUM'
$ws.Range('C40').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.789'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.26%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.490'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.09%  '
$ws.Range('E42').Value = '  -1.66%  '
$ws.Range('E43').Value = '  +0.75%  '
$ws.Range('D44').Value = '1.735.06'
$ws.Range('E44').Value = '  +2.32%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '92.92'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.23%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.51'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.01%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '53.49'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.89%  '
$ws.Range('E48').Value = '  -0.24%  '
$ws.Range('E49').Value = '  +0.86%  '
$ws.Range('E50').Value = '  +0.07%  '
$ws.Range('E51').Value = '  -0.03%  '
